# The workbook's "sheet1" has a block of repeating 4-column groups
# (Alain / Henri / Tony / Dulcinee headers in row 1, OUI/NON answers in
# rows 2-9) that runs from column E up to column ARP, immediately
# followed by an "Adresse de courriel" column and a trailing empty
# column.
#
# This edit inserts one more 4-column group repeated three times (i.e.
# 12 new columns) right before the email column, continuing the
# existing repeating pattern, and shifts the email / trailing-empty
# columns (and everything else) 12 columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Insert 12 blank columns at ARQ:ASB, pushing the email column
# (previously ARQ, now ASC) and everything after it to the right.
$ws.Range("ARQ1:ASB9").EntireColumn.Insert()

# Seed the newly inserted columns by copying the last existing
# 4-column cycle (ARM:ARP) - this carries over both the cell style
# and the Alain/Henri/Tony/Dulcinee (row 1) / OUI/NON (rows 2-9)
# values - and tiling it 3 times to fill all 12 new columns.
$src = $ws.Range("ARM1:ARP9")
$src.Copy()
$ws.Range("ARQ1:ART9").PasteSpecial(-4123)
$ws.Range("ARU1:ARX9").PasteSpecial(-4123)
$ws.Range("ARY1:ASB9").PasteSpecial(-4123)
$excel.CutCopyMode = 0
